# "Problem the result writing fixed" - rewrite the accuracy-score results
# table: the run's metadata (language/issue date/iteration count) changes,
# the BertTransformer row gets a real "51.5*" score, the Word2Vec row is
# replaced by a second (now-blank) BertTransformer run, and the
# Doc2Vec/TF-unigram rows disappear entirely since only two runs were made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- General info: issue date moved forward a day ---
$ws.Range("A5").Value = "Issue date: 22/12/2020 10:37:21"

# --- Row 16: first BertTransformer run, now in Hebrew, 20 iterations ---
$ws.Range("B16").Value = "Hebrew"
$ws.Range("C16").Value = "Stylistic Features: acf,  aof,  caf,  e50th,  fdf,  frc,  huf,  mef,  vof"
$ws.Range("D16").Value = "None"
$ws.Range("E16").Value = "lowercase"
$ws.Range("F16").Value = "5 folds X 20 iterations CV"

# The old "nan" placeholder in J16 moves to H16 and becomes a real score,
# styled like the other red significance-marked results (style of F11).
$ws.Range("J16").Clear()
$ws.Range("F11").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = "51.5*"

# --- Row 17: replaces the old Word2VecTransfomer row with a second
#     BertTransformer run (same settings, still unscored) ---
$ws.Range("H17:J17").Clear()
$ws.Range("B17").Value = "Hebrew"
$ws.Range("C17").Value = "BertTransformer"
$ws.Range("D17").Value = "None"
$ws.Range("E17").Value = "lowercase"
$ws.Range("F17").Value = "5 folds X 20 iterations CV"

# G17 already carries the blue significance-marker style from the old
# Word2Vec row (70.1*), which is exactly the style this "nan" needs.
$ws.Range("G17").Value = "nan"

# --- Rows 18 (Doc2Vec/TF-unigram extra runs) & 19 no longer exist ---
$ws.Rows("18:19").Delete()

# --- Column C needs to be much wider to fit the long stylistic-features text ---
$ws.Columns("C").ColumnWidth = 72.75

# --- Table style refresh ---
$tbl = $ws.ListObjects.Item(1)
$tbl.TableStyle = "TableStyleLight10"
